$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7; existing rows 7.. shift down to 8..
$ws.Rows("7:7").Insert()

# Populate the newly inserted row 7 with the new record's data.
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C7").Value = "Arica y Parinacota"
$ws.Range("D7").Value = 45282
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = 100112027
$ws.Range("G7").Value = "Melón"
$ws.Range("H7").Value = "Tuna"
$ws.Range("I7").Value = "Cuarta"
$ws.Range("J7").Value = 170
$ws.Range("K7").Value = 10000
$ws.Range("L7").Value = 12000
$ws.Range("M7").Value = 11059
$ws.Range("N7").Value = "`$/caja 40 unidades"
$ws.Range("O7").Value = "Región de Arica y Parinacota"
$ws.Range("P7").Value = 276
$ws.Range("Q7").Value = 40
$ws.Range("R7").Value = "Hortaliza"
